$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = 45992
$ws.Cells.Item(2, 2).Value = 11006.6321200206
$ws.Cells.Item(2, 3).Value = 10182.8165762619
$ws.Cells.Item(2, 4).Value = 17472.26
$ws.Cells.Item(2, 5).Value = 6931.47273134514
$ws.Cells.Item(2, 6).Value = -14.9154455163733

# Row 3
$ws.Cells.Item(3, 1).Value = 45993
$ws.Cells.Item(3, 2).Value = 11025.4261617451
$ws.Cells.Item(3, 3).Value = 10167.0099993937
$ws.Cells.Item(3, 4).Value = 11232.26
$ws.Cells.Item(3, 5).Value = 6977.20654250875
$ws.Cells.Item(3, 6).Value = 246.331522579267

# Row 4
$ws.Cells.Item(4, 1).Value = 45994
$ws.Cells.Item(4, 2).Value = 11104.2395902848
$ws.Cells.Item(4, 3).Value = 10227.7393428217
$ws.Cells.Item(4, 4).Value = 11232.26
$ws.Cells.Item(4, 5).Value = 7041.58526410658
$ws.Cells.Item(4, 6).Value = 251.544358622011

# Row 5
$ws.Cells.Item(5, 1).Value = 45995
$ws.Cells.Item(5, 2).Value = 10950.9995684019
$ws.Cells.Item(5, 3).Value = 10066.8342310277
$ws.Cells.Item(5, 4).Value = 11232.26
$ws.Cells.Item(5, 5).Value = 6931.34546915175
$ws.Cells.Item(5, 6).Value = 240.246654174142

# Row 6
$ws.Cells.Item(6, 1).Value = 45996
$ws.Cells.Item(6, 2).Value = 10752.9511990924
$ws.Cells.Item(6, 3).Value = 9384.46339333456
$ws.Cells.Item(6, 4).Value = 11232.26
$ws.Cells.Item(6, 5).Value = 6783.90799927249
$ws.Cells.Item(6, 6).Value = 205.671308025294

# Row 7
$ws.Cells.Item(7, 1).Value = 45997
$ws.Cells.Item(7, 2).Value = 8391.01684689587
$ws.Cells.Item(7, 3).Value = 8918.90133105568
$ws.Cells.Item(7, 4).Value = 11232.26
$ws.Cells.Item(7, 5).Value = 7098.75248917023
$ws.Cells.Item(7, 6).Value = 199.391409176079

# Row 8
$ws.Cells.Item(8, 1).Value = 45998
$ws.Cells.Item(8, 2).Value = 8291.57854635235
$ws.Cells.Item(8, 3).Value = 8843.49263966442
$ws.Cells.Item(8, 4).Value = 11232.26
$ws.Cells.Item(8, 5).Value = 7083.68516255095
$ws.Cells.Item(8, 6).Value = 195.621575092307

# Row 9
$ws.Cells.Item(9, 1).Value = 45999
$ws.Cells.Item(9, 2).Value = 8451.27859756152
$ws.Cells.Item(9, 3).Value = 8955.72420074943
$ws.Cells.Item(9, 4).Value = 11232.26
$ws.Cells.Item(9, 5).Value = 7517.40805337532
$ws.Cells.Item(9, 6).Value = 218.369677255198

# Row 10
$ws.Cells.Item(10, 1).Value = 46000
$ws.Cells.Item(10, 2).Value = 9794.31262154814
$ws.Cells.Item(10, 3).Value = 9784.80147188389
$ws.Cells.Item(10, 4).Value = 11232.26
$ws.Cells.Item(10, 5).Value = 7941.32071217855
$ws.Cells.Item(10, 6).Value = 270.577591002601

# Row 11
$ws.Cells.Item(11, 1).Value = 46001
$ws.Cells.Item(11, 2).Value = 9794.31262154814
$ws.Cells.Item(11, 3).Value = 9841.04134706554
$ws.Cells.Item(11, 4).Value = 11232.26
$ws.Cells.Item(11, 5).Value = 7941.32071217855
$ws.Cells.Item(11, 6).Value = 272.92091913517

# Row 12
$ws.Cells.Item(12, 1).Value = 46002
$ws.Cells.Item(12, 2).Value = 9794.31262154814
$ws.Cells.Item(12, 3).Value = 9730.40843290942
$ws.Cells.Item(12, 4).Value = 11232.26
$ws.Cells.Item(12, 5).Value = 7941.32071217855
$ws.Cells.Item(12, 6).Value = 268.311214378665

# Row 13
$ws.Cells.Item(13, 1).Value = 46003
$ws.Cells.Item(13, 2).Value = 9794.31262154814
$ws.Cells.Item(13, 3).Value = 8822.34103794774
$ws.Cells.Item(13, 4).Value = 11232.26
$ws.Cells.Item(13, 5).Value = 7941.32071217855
$ws.Cells.Item(13, 6).Value = 230.475072921929

# Row 14
$ws.Cells.Item(14, 1).Value = 46004
$ws.Cells.Item(14, 2).Value = 8553.43567882279
$ws.Cells.Item(14, 3).Value = 8814.96952263271
$ws.Cells.Item(14, 4).Value = 11232.26
$ws.Cells.Item(14, 5).Value = 7532.8353075375
$ws.Cells.Item(14, 6).Value = 213.147701257092

# Row 15
$ws.Cells.Item(15, 1).Value = 46005
$ws.Cells.Item(15, 2).Value = 8451.27859756152
$ws.Cells.Item(15, 3).Value = 9333.23912320702
$ws.Cells.Item(15, 4).Value = 11232.26
$ws.Cells.Item(15, 5).Value = 7517.48465969305
$ws.Cells.Item(15, 6).Value = 234.102657620836
